# Weekly update: insert two new data rows (Primera / Segunda) for the
# latest market date, pushing the existing historical rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 262 - everything
# from (old) row 262 downward shifts down by two rows.
$ws.Rows("262:263").Insert()

# --- New row 262: Betarraga, Primera, fecha 44813 ---
$ws.Range("A262").Value = 5
$ws.Range("B262").Value = "Macroferia Regional de Talca"
$ws.Range("C262").Value = "Maule"
$ws.Range("D262").Value = 44813
$ws.Range("E262").Value = 7
$ws.Range("F262").Value = 100114014
$ws.Range("G262").Value = "Betarraga"
$ws.Range("H262").Value = "Sin especificar"
$ws.Range("I262").Value = "Primera"
$ws.Range("J262").Value = 2000
$ws.Range("K262").Value = 1000
$ws.Range("L262").Value = 1000
$ws.Range("M262").Value = 1000
$ws.Range("N262").Value = "$/paquete 5 unidades"
$ws.Range("O262").Value = "Región del Maule"
$ws.Range("P262").Value = 200
$ws.Range("Q262").Value = 5
$ws.Range("R262").Value = "Hortaliza"

# --- New row 263: Betarraga, Segunda, fecha 44813 ---
$ws.Range("A263").Value = 5
$ws.Range("B263").Value = "Macroferia Regional de Talca"
$ws.Range("C263").Value = "Maule"
$ws.Range("D263").Value = 44813
$ws.Range("E263").Value = 7
$ws.Range("F263").Value = 100114014
$ws.Range("G263").Value = "Betarraga"
$ws.Range("H263").Value = "Sin especificar"
$ws.Range("I263").Value = "Segunda"
$ws.Range("J263").Value = 2000
$ws.Range("K263").Value = 800
$ws.Range("L263").Value = 800
$ws.Range("M263").Value = 800
$ws.Range("N263").Value = "$/paquete 5 unidades"
$ws.Range("O263").Value = "Región del Maule"
$ws.Range("P263").Value = 160
$ws.Range("Q263").Value = 5
$ws.Range("R263").Value = "Hortaliza"
